$p = $ppt.ActivePresentation

# --- Slide 1: "Area of Interest" bullet list ---
# Merge "Big Data and Data Platforms:" / "from collection to exploitation;"
# into a single bullet "Big Data and Data Platforms" (drop the sub-bullet).
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(2).TextFrame.TextRange
$tr1.Paragraphs(6, 1).Delete()
$tr1.Paragraphs(5, 1).Text = "Big Data and Data Platforms"

# --- Slide 3: "Future plans" bullet list ---
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange
$tr3.Paragraphs(2, 1).Runs(1, 1).Text = "Hopefully finish and publish the paper;"
$tr3.Paragraphs(3, 1).Runs(1, 1).Text = "optimize querying performance on the data structure."
